$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and G hold numeric-looking values that are actually stored as
# text in this sheet (e.g. "236.94", "15"). Force text format on the full
# data range first so Excel doesn't auto-convert the new values to numbers
# when they are assigned below.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Column G: every data row (2-51) moves from "15" to "16" (hour refresh).
$ws.Range("G2:G51").Value = "16"

# Column D: updated price snapshots for the rows that changed.
$ws.Range("D2").Value = "236.47"
$ws.Range("D3").Value = "21.91"
$ws.Range("D4").Value = "5.379"
$ws.Range("D5").Value = "0.05614"
$ws.Range("D6").Value = "6.475"
$ws.Range("D7").Value = "3.350"
$ws.Range("D8").Value = "0.7998"
$ws.Range("D9").Value = "1.040"
$ws.Range("D11").Value = "0.07296"
$ws.Range("D12").Value = "0.03120"
$ws.Range("D14").Value = "0.09241"
$ws.Range("D15").Value = "0.001665"
$ws.Range("D16").Value = "3.255"
$ws.Range("D17").Value = "0.04765"
$ws.Range("D19").Value = "0.006240"
$ws.Range("D20").Value = "0.005070"
$ws.Range("D21").Value = "0.001052"
$ws.Range("D23").Value = "0.0003902"
$ws.Range("D24").Value = "3.962"
$ws.Range("D40").Value = "0.04086"
$ws.Range("D41").Value = "0.007008"
$ws.Range("D44").Value = "0.008845"
$ws.Range("D47").Value = "0.6755"
$ws.Range("D48").Value = "0.03683"
$ws.Range("D50").Value = "0.01010"
